$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.432.86'
$ws.Range("E2").Value = '  -0.99%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.521.05'
$ws.Range("E3").Value = '  -0.78%  '
$ws.Range("E4").Value = '  -0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.42'
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.07'
$ws.Range("E6").Value = '  -3.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  +0.30%  '
$ws.Range("E8").Value = '  -1.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.526.13'
$ws.Range("E9").Value = '  -1.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.101'
$ws.Range("E10").Value = '  -0.13%  '
$ws.Range("E11").Value = '  +0.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.42'
$ws.Range("E12").Value = '  -2.25%  '
$ws.Range("E13").Value = '  -3.41%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.973.82'
$ws.Range("E14").Value = '  -0.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.34'
$ws.Range("E15").Value = '  -2.41%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '59.325.07'
$ws.Range("E16").Value = '  -1.00%  '
$ws.Range("E17").Value = '  -1.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.515.03'
$ws.Range("E18").Value = '  -1.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.10'
$ws.Range("E19").Value = '  -2.33%  '
$ws.Range("E20").Value = '  -1.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '325.45'
$ws.Range("E21").Value = '  -0.93%  '
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("E23").Value = '  -0.97%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.31'
$ws.Range("E24").Value = '  +0.91%  '
$ws.Range("E25").Value = '  -4.70%  '
$ws.Range("E27").Value = '  +1.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.73'
$ws.Range("E28").Value = '  -3.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0785'
$ws.Range("E29").Value = '  -2.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.72'
$ws.Range("E30").Value = '  -5.67%  '
$ws.Range("E31").Value = '  -1.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '165.37'
$ws.Range("E32").Value = '  +1.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.12'
$ws.Range("E34").Value = '  -9.64%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.41'
$ws.Range("E35").Value = '  -6.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.58'
$ws.Range("E36").Value = '  -1.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.17'
$ws.Range("E37").Value = '  -7.01%  '
$ws.Range("E38").Value = '  -2.79%  '
$ws.Range("E39").Value = '  -1.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.811'
$ws.Range("E40").Value = '  -3.73%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.24'
$ws.Range("E41").Value = '  -9.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '281.62'
$ws.Range("E42").Value = '  -7.33%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.997'
$ws.Range("E43").Value = '  +0.37%  '
$ws.Range("E44").Value = '  -1.45%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.89'
$ws.Range("E45").Value = '  +0.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '126.28'
$ws.Range("E46").Value = '  +0.96%  '
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("E48").Value = '  -1.83%  '
$ws.Range("E49").Value = '  -2.57%  '
$ws.Range("E50").Value = '  -2.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.772.36'
$ws.Range("E51").Value = '  -2.80%  '
